$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 600
$ws.Range("I46").Value = 400
$ws.Range("K46").Value = 1200
$ws.Range("M46").Value = -1081

$ws.Range("H60").Value = 600
$ws.Range("I60").Value = 400
$ws.Range("K60").Value = 1200
$ws.Range("M60").Value = -716

$ws.Range("H70").Value = 1384.2307
$ws.Range("I70").Value = 1124.375
$ws.Range("J70").Value = 1800
$ws.Range("K70").Value = 3373.125
$ws.Range("L70").Value = 5400
$ws.Range("M70").Value = -3103.125
$ws.Range("N70").Value = -5940

$ws.Range("H73").Value = 1384.2307
$ws.Range("I73").Value = 1124.375
$ws.Range("J73").Value = 1800
$ws.Range("K73").Value = 3373.125
$ws.Range("L73").Value = 5400
$ws.Range("M73").Value = -2437.125
$ws.Range("N73").Value = -7272

$ws.Range("H107").Value = 856.5
$ws.Range("I107").Value = 863
$ws.Range("J107").Value = 815.3333
$ws.Range("K107").Value = 863
$ws.Range("L107").Value = 815.3333
$ws.Range("M107").Value = 1057
$ws.Range("N107").Value = -4655.3333

$ws.Range("H111").Value = 4005.3333
$ws.Range("I111").Value = 6599.3335
$ws.Range("J111").Value = 2708.3333
$ws.Range("K111").Value = 19798.0005
$ws.Range("L111").Value = 8124.999899999999
$ws.Range("M111").Value = -16731.0005
$ws.Range("N111").Value = -14258.9999

$ws.Range("H113").Value = 200006160
$ws.Range("I113").Value = 500000900
$ws.Range("J113").Value = 9666.666999999999
$ws.Range("K113").Value = 500000900
$ws.Range("L113").Value = 9666.666999999999
$ws.Range("N113").Value = -16174.667
$ws.Range("M113").Value = -499997646

$ws.Range("H129").Value = 747.1
$ws.Range("I129").Value = 360
$ws.Range("J129").Value = 802.4
$ws.Range("K129").Value = 1080
$ws.Range("L129").Value = 2407.2
$ws.Range("M129").Value = 3920
$ws.Range("N129").Value = -12407.2

$ws.Range("H137").Value = 1652.1666
$ws.Range("I137").Value = 1768.875
$ws.Range("J137").Value = 1418.75
$ws.Range("K137").Value = 5306.625
$ws.Range("L137").Value = 4256.25
$ws.Range("M137").Value = -2756.625
$ws.Range("N137").Value = -9356.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 711
$ws.Range("I2").Value = 748.2727
$ws.Range("K2").Value = 748.2727
$ws.Range("M2").Value = -635.2727

$ws.Range("H32").Value = 5627.8374
$ws.Range("I32").Value = 4352.465
$ws.Range("J32").Value = 11664.6
$ws.Range("K32").Value = 4352.465
$ws.Range("L32").Value = 11664.6
$ws.Range("M32").Value = -4065.465
$ws.Range("N32").Value = -12238.6

$ws.Range("H45").Value = 3661.6538
$ws.Range("I45").Value = 3739.7
$ws.Range("J45").Value = 3612.875
$ws.Range("K45").Value = 3739.7
$ws.Range("L45").Value = 3612.875
$ws.Range("M45").Value = -3362.7
$ws.Range("N45").Value = -4366.875

$ws.Range("H63").Value = 1954671
$ws.Range("I63").Value = 1649.0667
$ws.Range("J63").Value = 31250000
$ws.Range("K63").Value = 1649.0667
$ws.Range("L63").Value = 31250000
$ws.Range("M63").Value = -963.0667000000001
$ws.Range("N63").Value = -31251372

$ws.Range("H66").Value = 1954671
$ws.Range("I66").Value = 1649.0667
$ws.Range("J66").Value = 31250000
$ws.Range("K66").Value = 8245.333500000001
$ws.Range("L66").Value = 156250000
$ws.Range("M66").Value = -4813.333500000001
$ws.Range("N66").Value = -156256864

$ws.Range("H116").Value = 711
$ws.Range("I116").Value = 748.2727
$ws.Range("K116").Value = 748.2727
$ws.Range("M116").Value = 1545.7273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 711
$ws.Range("I3").Value = 748.2727
$ws.Range("K3").Value = 748.2727
$ws.Range("M3").Value = -634.2727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1713.875
$ws.Range("I16").Value = 1451.8334
$ws.Range("K16").Value = 1451.8334
$ws.Range("M16").Value = -1164.8334

$ws.Range("H31").Value = 3989.8215
$ws.Range("I31").Value = 947.2727
$ws.Range("J31").Value = 5958.5293
$ws.Range("K31").Value = 947.2727
$ws.Range("L31").Value = 5958.5293
$ws.Range("M31").Value = -652.2727
$ws.Range("N31").Value = -6548.5293

$ws.Range("H34").Value = 3989.8215
$ws.Range("I34").Value = 947.2727
$ws.Range("J34").Value = 5958.5293
$ws.Range("K34").Value = 947.2727
$ws.Range("L34").Value = 5958.5293
$ws.Range("M34").Value = -745.2727
$ws.Range("N34").Value = -6362.5293

$ws.Range("H113").Value = 1713.875
$ws.Range("I113").Value = 1451.8334
$ws.Range("K113").Value = 1451.8334
$ws.Range("M113").Value = 718.1666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1106.2433
$ws.Range("I5").Value = 899.3333
$ws.Range("K5").Value = 2697.9999
$ws.Range("M5").Value = -2585.9999

$ws.Range("H33").Value = 67
$ws.Range("I33").Value = 4
$ws.Range("J33").Value = 79.59999999999999
$ws.Range("K33").Value = 24
$ws.Range("L33").Value = 477.6
$ws.Range("M33").Value = 259
$ws.Range("N33").Value = -1043.6

$ws.Range("H135").Value = 1106.2433
$ws.Range("I135").Value = 899.3333
$ws.Range("K135").Value = 8093.9997
$ws.Range("M135").Value = -5558.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4457.3335
$ws.Range("I113").Value = 6008.9
$ws.Range("J113").Value = 2070.3076
$ws.Range("K113").Value = 6008.9
$ws.Range("L113").Value = 2070.3076
$ws.Range("M113").Value = -3838.9
$ws.Range("N113").Value = -6410.3076

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2344.4443
$ws.Range("I22").Value = 2816.8333
$ws.Range("J22").Value = 1399.6666
$ws.Range("K22").Value = 2816.8333
$ws.Range("L22").Value = 1399.6666
$ws.Range("M22").Value = -2521.8333
$ws.Range("N22").Value = -1989.6666

$ws.Range("H27").Value = 2344.4443
$ws.Range("I27").Value = 2816.8333
$ws.Range("J27").Value = 1399.6666
$ws.Range("K27").Value = 2816.8333
$ws.Range("L27").Value = 1399.6666
$ws.Range("M27").Value = -2709.8333
$ws.Range("N27").Value = -1613.6666

$ws.Range("H46").Value = 1111.0209
$ws.Range("I46").Value = 1090.711
$ws.Range("J46").Value = 1415.6666
$ws.Range("K46").Value = 1090.711
$ws.Range("L46").Value = 1415.6666
$ws.Range("M46").Value = -902.711
$ws.Range("N46").Value = -1791.6666

$ws.Range("H132").Value = 2351.6316
$ws.Range("I132").Value = 1408.5834
$ws.Range("J132").Value = 3968.2856
$ws.Range("K132").Value = 4225.7502
$ws.Range("L132").Value = 11904.8568
$ws.Range("M132").Value = -1695.7502
$ws.Range("N132").Value = -16964.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1594.6957
$ws.Range("I122").Value = 1561.7894
$ws.Range("J122").Value = 1751
$ws.Range("K122").Value = 4685.3682
$ws.Range("L122").Value = 5253
$ws.Range("M122").Value = -2235.3682
$ws.Range("N122").Value = -10153

$ws.Range("H123").Value = 30286
$ws.Range("J123").Value = 30286
$ws.Range("L123").Value = 30286
$ws.Range("N123").Value = -40086

$ws.Range("H126").Value = 1525.1613
$ws.Range("I126").Value = 1145.4231
$ws.Range("K126").Value = 3436.2693
$ws.Range("M126").Value = -966.2692999999999
